$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both contain the same table of events in
# rows 2-5. Column F holds the "想去人数" (want-to-go count), which was
# bumped slightly for four rows on each of those sheets.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 5580
    $ws.Range("F3").Value = 171
    $ws.Range("F4").Value = 949
    $ws.Range("F5").Value = 16
}
